# Attempt VII follow-up: add an adjusted-speed section (columns G:H) to
# Sheet1. G1 holds the 10%-reduction multiplier (0.9); G3:G15 multiply the
# default speed (B) by the revised correction factor (E) and the G1
# multiplier; H3:H15 express that adjusted speed as a fraction of the
# original default speed (B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiplier used by column G (10% reduction factor).
$ws.Range("G1").Value = 0.9

# Column G: B * E * $G$1  (absolute reference to the multiplier in G1).
$ws.Range("G3").Formula = "=B3*E3*`$G`$1"
$ws.Range("G4:G15").Formula = "=B4*E4*`$G`$1"

# Column H: G / B  (adjusted speed as a fraction of the default speed).
$ws.Range("H3").Formula = "=G3/B3"
$ws.Range("H4:H15").Formula = "=G4/B4"

# Leave the selection where the author left it after entering the data.
$ws.Range("H10").Select() | Out-Null
